$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -6
$ws.Range("F9").Value = 6
$ws.Range("F12").Value = -5
$ws.Range("F16").Value = -1
